$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Cases query, B2): append "order By ... LIMIT 100" clause ---
$b2 = $ws.Range("B2").Value()
$ws.Range("B2").Value = $b2 + "`n order By ss.study_subject_id ASC LIMIT 100 "
$ws.Rows("2").RowHeight = 331.2

# --- Row 3 (Samples query, B3): append "order By ... LIMIT 100" clause ---
$b3 = $ws.Range("B3").Value()
$ws.Range("B3").Value = $b3 + "`n order By samp.sample_id ASC LIMIT 100"
$ws.Rows("3").RowHeight = 360

# --- Row 4 (Files query, B4): replace trailing "order by f.file_name" with the new clause ---
$b4 = $ws.Range("B4").Value()
$b4New = $b4.Replace("    order by f.file_name", "  order By f.file_name ASC LIMIT 100")
$ws.Range("B4").Value = $b4New
